$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Formula = "'245.63"
$ws.Range('D3').Formula = "'25.46"
$ws.Range('D4').Formula = "'5.129"
$ws.Range('D5').Formula = "'0.05591"
$ws.Range('D6').Formula = "'6.488"
$ws.Range('D7').Formula = "'3.027"
$ws.Range('D8').Formula = "'0.8173"
$ws.Range('D9').Formula = "'0.8505"
$ws.Range('B10').Value = 'One'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D10').Formula = "'0.009759"
$ws.Range('E10').Value = '9OneONEBestin24h'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Formula = "'0.1340"
$ws.Range('E11').Value = '10WazirXWRX'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').Formula = "'0.02849"
$ws.Range('E12').Value = '11BitrueCoinBTR'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').Formula = "'0.09405"
$ws.Range('E13').Value = '12BitMartTokenBMX'
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').Formula = "'0.001526"
$ws.Range('E14').Value = '13BitForexTokenBF'
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D15').Formula = "'0.006193"
$ws.Range('E15').Value = '14TigerCashTCH'
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').Formula = "'3.533"
$ws.Range('E16').Value = '15LEOLEO'
$ws.Range('B17').Value = 'BTSEToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D17').Formula = "'2.118"
$ws.Range('E17').Value = '16BTSETokenBTSE'
$ws.Range('B18').Value = 'BitpandaEcosystemToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D18').Formula = "'0.3168"
$ws.Range('E18').Value = '17BitpandaEcosystemTokenBEST'
$ws.Range('B19').Value = 'MandalaExchangeToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D19').Formula = "'0.06950"
$ws.Range('E19').Value = '18MandalaExchangeTokenMDX'
$ws.Range('B20').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C20').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D20').Formula = "'0.03234"
$ws.Range('E20').Value = '19LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('D22').Formula = "'3.755"
$ws.Range('D23').Formula = "'0.04692"
$ws.Range('D25').Formula = "'0.001247"
$ws.Range('D26').Formula = "'0.004606"
$ws.Range('D27').Formula = "'0.00009595"
$ws.Range('E27').Value = '26NitroExNTX'
$ws.Range('D40').Formula = "'0.03657"
$ws.Range('B41').Value = 'KickToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D41').Formula = "'0.006137"
$ws.Range('E41').Value = '40KickTokenKICK'
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D42').Formula = "'0.1053"
$ws.Range('E42').Value = '41BKEXTokenBKK'
$ws.Range('D44').Formula = "'0.007389"
$ws.Range('D45').Formula = "'0.00005309"
$ws.Range('E47').Value = '46CoinbaseStockTokenCOINWorstin24h'
$ws.Range('D48').Formula = "'0.002126"
